$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2024-08-26 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-08-27 Tuesday", 2) | Out-Null

# Update each table cell by position (handles duplicate values safely)
$t = $d.Tables(1)
$t.Cell(1,1).Range.Text = "39-1=38"
$t.Cell(1,2).Range.Text = "69+15=84"
$t.Cell(1,3).Range.Text = "86-1=85"
$t.Cell(1,4).Range.Text = "7+82=89"
$t.Cell(1,5).Range.Text = "80-26=54"
$t.Cell(2,1).Range.Text = "94-0=94"
$t.Cell(2,2).Range.Text = "3+67=70"
$t.Cell(2,3).Range.Text = "6+44=50"
$t.Cell(2,4).Range.Text = "86-25=61"
$t.Cell(2,5).Range.Text = "89-17=72"
$t.Cell(3,1).Range.Text = "99-80=19"
$t.Cell(3,2).Range.Text = "15+77=92"
$t.Cell(3,3).Range.Text = "64-29=35"
$t.Cell(3,4).Range.Text = "45-3=42"
$t.Cell(3,5).Range.Text = "61-34=27"
$t.Cell(4,1).Range.Text = "53-17=36"
$t.Cell(4,2).Range.Text = "96-41=55"
$t.Cell(4,3).Range.Text = "20+4=24"
$t.Cell(4,4).Range.Text = "73+2=75"
$t.Cell(4,5).Range.Text = "25+66=91"
$t.Cell(5,1).Range.Text = "44+54=98"
$t.Cell(5,2).Range.Text = "62-12=50"
$t.Cell(5,3).Range.Text = "54-27=27"
$t.Cell(5,4).Range.Text = "39+24=63"
$t.Cell(5,5).Range.Text = "59-0=59"
$t.Cell(6,1).Range.Text = "90+1=91"
$t.Cell(6,2).Range.Text = "97-17=80"
$t.Cell(6,3).Range.Text = "66-60=6"
$t.Cell(6,4).Range.Text = "28-2=26"
$t.Cell(6,5).Range.Text = "18+20=38"
$t.Cell(7,1).Range.Text = "29+60=89"
$t.Cell(7,2).Range.Text = "81-51=30"
$t.Cell(7,3).Range.Text = "60+26=86"
$t.Cell(7,4).Range.Text = "37+30=67"
$t.Cell(7,5).Range.Text = "28-10=18"
$t.Cell(8,1).Range.Text = "97-11=86"
$t.Cell(8,2).Range.Text = "56+12=68"
$t.Cell(8,3).Range.Text = "86-4=82"
$t.Cell(8,4).Range.Text = "97-0=97"
$t.Cell(8,5).Range.Text = "96-34=62"
$t.Cell(9,1).Range.Text = "2+36=38"
$t.Cell(9,2).Range.Text = "90-4=86"
$t.Cell(9,3).Range.Text = "57+23=80"
$t.Cell(9,4).Range.Text = "61+38=99"
$t.Cell(9,5).Range.Text = "85-11=74"
$t.Cell(10,1).Range.Text = "1+68=69"
$t.Cell(10,2).Range.Text = "0+84=84"
$t.Cell(10,3).Range.Text = "87-82=5"
$t.Cell(10,4).Range.Text = "42-13=29"
$t.Cell(10,5).Range.Text = "42+3=45"
$t.Cell(11,1).Range.Text = "64+13=77"
$t.Cell(11,2).Range.Text = "4+9=13"
$t.Cell(11,3).Range.Text = "10+71=81"
$t.Cell(11,4).Range.Text = "11+65=76"
$t.Cell(11,5).Range.Text = "74-33=41"
$t.Cell(12,1).Range.Text = "84-5=79"
$t.Cell(12,2).Range.Text = "31+50=81"
$t.Cell(12,3).Range.Text = "42+25=67"
$t.Cell(12,4).Range.Text = "30+50=80"
$t.Cell(12,5).Range.Text = "89-60=29"
$t.Cell(13,1).Range.Text = "37+20=57"
$t.Cell(13,2).Range.Text = "2+39=41"
$t.Cell(13,3).Range.Text = "71-18=53"
$t.Cell(13,4).Range.Text = "17+44=61"
$t.Cell(13,5).Range.Text = "19+56=75"
$t.Cell(14,1).Range.Text = "73+0=73"
$t.Cell(14,2).Range.Text = "51+6=57"
$t.Cell(14,3).Range.Text = "89-1=88"
$t.Cell(14,4).Range.Text = "62-22=40"
$t.Cell(14,5).Range.Text = "66-10=56"
$t.Cell(15,1).Range.Text = "65-53=12"
$t.Cell(15,2).Range.Text = "21+71=92"
$t.Cell(15,3).Range.Text = "40+3=43"
$t.Cell(15,4).Range.Text = "40+18=58"
$t.Cell(15,5).Range.Text = "29-13=16"
$t.Cell(16,1).Range.Text = "30-7=23"
$t.Cell(16,2).Range.Text = "6-1=5"
$t.Cell(16,3).Range.Text = "97-59=38"
$t.Cell(16,4).Range.Text = "8+56=64"
$t.Cell(16,5).Range.Text = "61-25=36"
$t.Cell(17,1).Range.Text = "9+63=72"
$t.Cell(17,2).Range.Text = "72+20=92"
$t.Cell(17,3).Range.Text = "74-57=17"
$t.Cell(17,4).Range.Text = "14+53=67"
$t.Cell(17,5).Range.Text = "76+13=89"
$t.Cell(18,1).Range.Text = "20-18=2"
$t.Cell(18,2).Range.Text = "57-22=35"
$t.Cell(18,3).Range.Text = "79-70=9"
$t.Cell(18,4).Range.Text = "97-79=18"
$t.Cell(18,5).Range.Text = "62-53=9"
$t.Cell(19,1).Range.Text = "79-25=54"
$t.Cell(19,2).Range.Text = "49-8=41"
$t.Cell(19,3).Range.Text = "9+20=29"
$t.Cell(19,4).Range.Text = "68-48=20"
$t.Cell(19,5).Range.Text = "49-42=7"
$t.Cell(20,1).Range.Text = "42+35=77"
$t.Cell(20,2).Range.Text = "13+69=82"
$t.Cell(20,3).Range.Text = "53-12=41"
$t.Cell(20,4).Range.Text = "2+33=35"
$t.Cell(20,5).Range.Text = "73-52=21"
